$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$range = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
try {
  $tbl.TableStyle = "TableStyleMedium0"
  Write-Host "set ok"
} catch {
  Write-Host "ERR: $_"
}
